$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Edit 1: "Contato com o cliente deve ser feito em horário comercial,
# primeiramente por e-mail." -> remove ", primeiramente por e-mail" and
# split the trailing period into its own run (same formatting).
# --------------------------------------------------------------------
$d.Content.Find.Execute(", primeiramente por e-mail", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$pContato = $d.Paragraphs.Item(3).Range
$periodStart = $pContato.End - 2
$rPeriod = $d.Range($periodStart, $pContato.End - 1)
# Toggle a character property on/off (same net formatting) to force Word
# to materialize this as its own run without altering its appearance.
$rPeriod.Bold = 1
$rPeriod.Bold = 0

# --------------------------------------------------------------------
# Edit 2: remove the whole paragraph "Somente entrar em contato aos
# finais de semana caso seja algo de muita criticidade." entirely.
# --------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Delete()

# --------------------------------------------------------------------
# Edit 3: collapse the three runs of the "Caso haja necessidade..."
# paragraph (including the spell-check-flagged "Whatsapp" run) into a
# single run with the full sentence, dropping the proofErr markers.
# --------------------------------------------------------------------
$pCaso = $d.Paragraphs.Item(4).Range
$paraStart = $pCaso.Start
$paraEnd = $pCaso.End - 1

$firstPart = $d.Range($paraStart, $pCaso.End)
$firstPart.Find.Execute("Caso haja necessidade e falta de resposta via e-mail,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1End = $firstPart.End

$d.Range($run1End, $paraEnd).Delete()
$d.Range($paraStart, $run1End).InsertAfter(" Whatsapp está liberado para contato. ")
